# Generate Report for Handback
#
# Updates the localization-status workbook to reflect a completed
# handback for the "5a04021b-...md" source file:
#   - Status changes from "Ready for handoff" to
#     "Handed back: in sync with en-US" (zh-cn and de-de sheets).
#   - "Latest Target File" / "Latest Handback File" / "Latest Handback
#     DateTime" columns are populated for rows 2 & 3 on both language
#     sheets, including a new hyperlink in the "Latest Target File"
#     column.
#   - The widened columns (Status col on Overview/zh-cn/de-de, and the
#     "Latest Target File" / "Latest Handback File" columns) are resized
#     to fit the new, longer content.

$wb = $excel.ActiveWorkbook

$hyperlinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c369923b4c499ce48bed20f49f6f5c30168884f9/e2e/5a04021b-a533-40d8-9da6-7aaf68baea15.md"
$hyperlinkDisplay = "5a04021b-a533-40d8-9da6-7aaf68baea15.md"

$statusText = "Handed back: in sync with en-US"

# ColumnWidth values below are chosen so the stored (raw OOXML) column
# width lands as close as possible to the authored target width; this
# runtime quantizes ColumnWidth to 1/6-character steps (offset by 5/6),
# so exact legacy fractional widths can't always be reproduced bit for
# bit, but this gets within a fraction of a character.
$width30 = 29.166666666666668   # -> stored width ~30     (target 29.9777047293527)
$width40 = 39.166666666666664   # -> stored width  40     (target 40)

# ---------------------------------------------------------------------
# Overview sheet: widen the zh-cn / de-de status columns (E, F)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $width30
$wsOverview.Columns.Item(6).ColumnWidth = $width30

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Widen Status (C), Latest Target File (I) and Latest Handback File (J)
$wsZh.Columns.Item(3).ColumnWidth = $width30
$wsZh.Columns.Item(9).ColumnWidth = $width40
$wsZh.Columns.Item(10).ColumnWidth = $width40

# Row 2 & 3: mark as handed back
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

# Latest Target File (I): hyperlink to the source md file
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $hyperlinkUrl, "", "", $hyperlinkDisplay)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $hyperlinkUrl, "", "", $hyperlinkDisplay)

# Latest Handback File (J): the generated xliff handback file
$wsZh.Range("J2").Value = "5a04021b-a533-40d8-9da6-7aaf68baea15.70ceb74897fb9812e6766e5e58af29a09d2a00f7.zh-cn.xlf"
$wsZh.Range("J3").Value = "5a04021b-a533-40d8-9da6-7aaf68baea15.70ceb74897fb9812e6766e5e58af29a09d2a00f7.zh-cn.xlf"

# Latest Handback DateTime (K)
$wsZh.Range("K2").Value = "2016-08-25 00:59:27"
$wsZh.Range("K3").Value = "2016-08-25 00:59:27"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Widen Status (C), Latest Target File (I) and Latest Handback File (J)
$wsDe.Columns.Item(3).ColumnWidth = $width30
$wsDe.Columns.Item(9).ColumnWidth = $width40
$wsDe.Columns.Item(10).ColumnWidth = $width40

# Row 2 & 3: mark as handed back
$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

# Latest Target File (I): hyperlink to the source md file
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $hyperlinkUrl, "", "", $hyperlinkDisplay)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $hyperlinkUrl, "", "", $hyperlinkDisplay)

# Latest Handback File (J): the generated xliff handback file
$wsDe.Range("J2").Value = "5a04021b-a533-40d8-9da6-7aaf68baea15.70ceb74897fb9812e6766e5e58af29a09d2a00f7.de-de.xlf"
$wsDe.Range("J3").Value = "5a04021b-a533-40d8-9da6-7aaf68baea15.70ceb74897fb9812e6766e5e58af29a09d2a00f7.de-de.xlf"

# Latest Handback DateTime (K)
$wsDe.Range("K2").Value = "2016-08-25 00:59:34"
$wsDe.Range("K3").Value = "2016-08-25 00:59:34"
